# The underlying save-data generator was changed to compute the "K" column
# (header in G1 = "K", formerly derived from a "Strike#" style counter) using
# a different metric (std/mean based s_vals calc). This regenerates the
# stored K values for every existing data row (rows 2-30) on the active
# sheet, leaving all other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 3
    7  = 3
    8  = 2
    9  = 2
    10 = 11
    11 = 4
    12 = 5
    13 = 1
    14 = 3
    15 = 5
    16 = 2
    17 = 5
    18 = 0
    19 = 6
    20 = 0
    21 = 2
    22 = 4
    23 = 3
    24 = 2
    25 = 2
    26 = 6
    27 = 2
    28 = 2
    29 = 1
    30 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
